# Added Configurable zero_before_threshold parameter to enable setting dims
# before noise_threshold or First Rise Point to 0.
#
# This recomputes, for each "Step3_DataPts_*" sheet, the
# First_Noticeable_Increase_Index (C), First_Noticeable_Increase_Cumulative_Value (E)
# and Pulse_Width (G) columns for rows 3-6 (signal segments 2-5) to reflect the
# new zero_before_threshold behavior.

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "Step3_DataPts_0.5",
    "Step3_DataPts_0.7",
    "Step3_DataPts_0.8",
    "Step3_DataPts_0.9"
)

# New First_Noticeable_Increase_Index (column C) values per row (3-6), same
# across all four threshold sheets.
$newC = @{
    3 = 45
    4 = 44
    5 = 44
    6 = 48
}

# New First_Noticeable_Increase_Cumulative_Value (column E) values per row (3-6),
# same across all four threshold sheets.
$newE = @{
    3 = 0.001255529070589158
    4 = 0.001802049215885368
    5 = 0.0008840628832351372
    6 = 0.03332298524579122
}

# New Pulse_Width (column G) values per row (3-6), per sheet (Point_Exceeds_Index
# in column D is unchanged, so G = D - new C).
$newG = @{
    "Step3_DataPts_0.5" = @{ 3 = 49; 4 = 51; 5 = 50; 6 = 46 }
    "Step3_DataPts_0.7" = @{ 3 = 63; 4 = 64; 5 = 65; 6 = 60 }
    "Step3_DataPts_0.8" = @{ 3 = 75; 4 = 80; 5 = 80; 6 = 75 }
    "Step3_DataPts_0.9" = @{ 3 = 113; 4 = 117; 5 = 118; 6 = 110 }
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in 3..6) {
        $ws.Cells.Item($row, 3).Value = $newC[$row]
        $ws.Cells.Item($row, 5).Value = $newE[$row]
        $ws.Cells.Item($row, 7).Value = $newG[$sheetName][$row]
    }
}
